# PCM_model 2D improvement: raise the D34 coefficient (Cliquido-phase slope
# driver used by the B25:B... running-sum formulas) from 60000 to 70000.
# This is the single user edit; every other changed cell in the workbook is
# Excel recalculating formulas (and chart caches) that depend on D34.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dati")
$ws.Activate()

# The actual edit: D34 60000 -> 70000
$ws.Range("D34").Value = 70000

# Recalculate so every dependent formula (B25:B66 and the charts) picks up
# the new value.
$excel.Calculate()

# Reflect where the author's cursor/viewport ended up after making the edit:
# scrolled up a bit and left the freshly-edited D34 cell selected.
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D34").Select()
